$p = $ppt.ActivePresentation

# Delete the 8 slides that are being replaced (positions 65-72).
# After deletion, the 5 slides that were at positions 73-77 shift up
# to become positions 65-69, matching the target sldIdLst exactly.
for ($i = 0; $i -lt 8; $i++) {
    $p.Slides.Item(65).Delete()
}
